$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'28.737.32"
$ws.Cells.Item(2, 5).Value = "  +1.85%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "'1.575.29"
$ws.Cells.Item(3, 5).Value = "  -0.68%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.12%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'213.54"
$ws.Cells.Item(5, 5).Value = "  +0.12%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'0.491"
$ws.Cells.Item(6, 5).Value = "  +0.19%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.16%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'44.76"
$ws.Cells.Item(8, 5).Value = "  +2.03%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'24.14"
$ws.Cells.Item(9, 5).Value = "  +1.10%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -1.07%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.57%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'0.0891"
$ws.Cells.Item(12, 5).Value = "  +0.42%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'1.800.19"
$ws.Cells.Item(13, 5).Value = "  -0.70%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'1.575.15"
$ws.Cells.Item(14, 5).Value = "  -0.65%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  -1.07%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'28.728.77"
$ws.Cells.Item(16, 5).Value = "  +1.78%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  -1.60%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'62.41"
$ws.Cells.Item(18, 5).Value = "  -1.19%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'231.19"
$ws.Cells.Item(19, 5).Value = "  +1.87%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -0.95%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -1.67%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  +0.04%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -4.45%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -1.20%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'2.04"
$ws.Cells.Item(25, 5).Value = "  +4.67%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'151.81"
$ws.Cells.Item(26, 5).Value = "  -0.06%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -0.63%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -1.27%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -2.10%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +0.07%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'0.0481"
$ws.Cells.Item(31, 5).Value = "  +2.24%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -1.72%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'3.22"
$ws.Cells.Item(33, 5).Value = "  -0.52%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -0.98%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'1.397.86"
$ws.Cells.Item(35, 5).Value = "  -0.15%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +2.48%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'1.54"
$ws.Cells.Item(37, 5).Value = "  -3.07%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +0.66%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +3.05%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -0.16%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.524"
$ws.Cells.Item(41, 5).Value = "  -2.98%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  +0.08%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "RenderToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(43, 4).Value = "'1.90"
$ws.Cells.Item(43, 5).Value = "  +1.69%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "ARBITRUM"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(44, 4).Value = "'0.794"
$ws.Cells.Item(44, 5).Value = "  -1.98%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'0.0468"
$ws.Cells.Item(45, 5).Value = "  +1.55%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -1.72%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'0.961"
$ws.Cells.Item(47, 5).Value = "  -1.97%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'63.30"
$ws.Cells.Item(48, 5).Value = "  -1.39%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'1.712.23"
$ws.Cells.Item(49, 5).Value = "  -0.57%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'86.55"
$ws.Cells.Item(50, 5).Value = "  -0.44%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'0.0₆0103"
$ws.Cells.Item(51, 5).Value = "  +0.09%  "

Write-Host "Update complete"